# Remove the "42" and "EV-1" bus route entries (original data rows 2-7),
# shifting the "61" bus route entries up so the sheet starts with bus 61.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:E7").EntireRow.Delete() | Out-Null
